$d = $word.ActiveDocument
$apos = [char]8217

# --- Change 1: append a new sentence to the end of paragraph 1 ---
# (keeps the same sz/szCs formatting as the rest of the paragraph by
# cloning the formatting of an existing same-style run before typing
# the new text into it)
$p1 = $d.Paragraphs(1)
$p1EndBeforeMark = $p1.Range.End - 1
$srcFmt = $d.Range($p1EndBeforeMark - 1, $p1EndBeforeMark)
$insertPoint = $d.Range($p1EndBeforeMark, $p1EndBeforeMark)
$insertPoint.FormattedText = $srcFmt.FormattedText
$newRun = $d.Range($p1EndBeforeMark, $p1EndBeforeMark + 1)
$newRun.Text = "I wanted her to hurt from loneliness as deeply as I had been hurt"

# --- Change 2: rework paragraph 2 ---

# 2a. Insert a new sentence after "...how I loved." and before "At the
# same time, ..."
$found1 = $d.Content.Find.Execute(
    "how I loved. At the same time",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "how I loved. I was curious about this God who loved so perfectly, so I went to services and family groups to learn more than anything else. At the same time",
    2)

# 2b. Replace the closing lines of the paragraph with the new, longer
# ending.
$replacement2 = "I wanted to go to morning prayer just to pray for my friends. " + `
    "I thought I was being as loving of a person as I could be. Then, I remembered Ivy, and I knew that Go" + `
    "d would want me to forgive her, and that was something that I couldn" + $apos + "t do. My love couldn" + $apos + `
    "t reach that far. God pursued me even more through class and accountability and discipleship. Through God" + $apos + "s pursuit, I realized that my love was a selfish love, that it was a love to fill a hole that only God could fill. When God began to fill that hole, His gospel became more and more real. " + `
    "I thought I couldn" + $apos + "t forgive Ivy, but after experiencing God" + $apos + "s love, I could finally forgive her. I thought I could love people my way, but people are limited and God is eternal"

$found2 = $d.Content.Find.Execute(
    "I wanted to go to morning prayer just to pray for my friends. I listened to a lot of sermons too. I thought I was doing all that I could to love and be loved.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    $replacement2,
    2)

# 2c. The paragraph used to end with a lone space run after the
# _GoBack bookmark; the new text now runs right up against the
# bookmark, so drop that trailing space.
$p2 = $d.Paragraphs(2)
$trailing = $d.Range($p2.Range.End - 2, $p2.Range.End - 1)
if ($trailing.Text -eq " ") {
    $trailing.Delete()
}

Write-Host "Find1:" $found1 "Find2:" $found2
